$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values in column D are plain numeric-looking strings (e.g. "325.00").
# Excel auto-converts those to real numbers on assignment (dropping the trailing
# zero / exact text), so force those specific cells to Text format first so the
# literal string is preserved, matching the source data feeds text formatting.

$ws.Range("D2").Value = "43.541.29"
$ws.Range("E2").Value = "  +2.70%  "
$ws.Range("D3").Value = "2.417.86"
$ws.Range("E3").Value = "  +8.93%  "
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.00"
$ws.Range("E5").Value = "  +12.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.99"
$ws.Range("E6").Value = "  -3.69%  "
$ws.Range("E7").Value = "  +4.10%  "
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.662"
$ws.Range("E9").Value = "  +11.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.32"
$ws.Range("E10").Value = "  -1.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0952"
$ws.Range("E11").Value = "  +4.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.70"
$ws.Range("E12").Value = "  +2.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.05"
$ws.Range("E13").Value = "  +3.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.22"
$ws.Range("E14").Value = "  +16.23%  "
$ws.Range("E15").Value = "  +2.84%  "
$ws.Range("D16").Value = "2.782.35"
$ws.Range("E16").Value = "  +8.90%  "
$ws.Range("D17").Value = "2.419.55"
$ws.Range("E17").Value = "  +8.42%  "
$ws.Range("D18").Value = "43.565.42"
$ws.Range("E18").Value = "  +2.97%  "
$ws.Range("E19").Value = "  +6.41%  "
$ws.Range("E20").Value = "  +5.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "75.65"
$ws.Range("E21").Value = "  +4.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.50"
$ws.Range("E22").Value = "  +4.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "260.75"
$ws.Range("E23").Value = "  +13.89%  "
$ws.Range("E24").Value = "  +4.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.64"
$ws.Range("E25").Value = "  +8.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.99"
$ws.Range("E26").Value = "  +5.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.98"
$ws.Range("E28").Value = "  +0.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.94"
$ws.Range("E29").Value = "  +10.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "180.11"
$ws.Range("E30").Value = "  +4.29%  "
$ws.Range("E31").Value = "  +2.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "38.36"
$ws.Range("E32").Value = "  +4.18%  "
$ws.Range("E33").Value = "  +2.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0936"
$ws.Range("E34").Value = "  +7.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.94"
$ws.Range("E35").Value = "  +7.11%  "
$ws.Range("E36").Value = "  +6.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.90"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0373"
$ws.Range("E38").Value = "  +2.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.02"
$ws.Range("E39").Value = "  -1.93%  "
$ws.Range("E40").Value = "  +3.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.90"
$ws.Range("E41").Value = "  +23.30%  "
$ws.Range("E42").Value = "  +27.24%  "
$ws.Range("E43").Value = "  +3.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "125.89"
$ws.Range("E44").Value = "  +25.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "69.83"
$ws.Range("E45").Value = "  -4.65%  "
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.63"
$ws.Range("E47").Value = "  +3.68%  "
$ws.Range("E48").Value = "  +7.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.55"
$ws.Range("E49").Value = "  +13.96%  "
$ws.Range("E50").Value = "  +4.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "75.00"
$ws.Range("E51").Value = "  +12.20%  "
